# The underlying data rows (sheet rows 2-12) are rotated: the two
# top records (old rows 2-3) move to the bottom of the block (new rows
# 11-12), and the remaining records (old rows 4-12) shift up to become
# new rows 2-10. This script captures the whole A2:AY12 block as a
# single 2-D array, rebuilds it in the new row order, and writes it
# back in one shot so that every column moves together with its row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRange = $ws.Range("A2:AY12")
$dstRange = $ws.Range("A2:AY12")

$arr = $srcRange.Value()

$rowCount = $arr.GetUpperBound(0)   # 11 (array rows for sheet rows 2..12)
$colCount = $arr.GetUpperBound(1)   # 51 (columns A..AY)

# Mapping of 1-based array row ($r, where $r=1 is sheet row 2) to the
# 1-based array row that supplies its new content, i.e. the rotation
# seen in the target workbook:
#   sheet row 2  ($r=1)  <- old sheet row 4  ($rowMap[1]=3)
#   sheet row 3  ($r=2)  <- old sheet row 5  ($rowMap[2]=4)
#   sheet row 4  ($r=3)  <- old sheet row 6  ($rowMap[3]=5)
#   sheet row 5  ($r=4)  <- old sheet row 7  ($rowMap[4]=6)
#   sheet row 6  ($r=5)  <- old sheet row 8  ($rowMap[5]=7)
#   sheet row 7  ($r=6)  <- old sheet row 9  ($rowMap[6]=8)
#   sheet row 8  ($r=7)  <- old sheet row 10 ($rowMap[7]=9)
#   sheet row 9  ($r=8)  <- old sheet row 11 ($rowMap[8]=10)
#   sheet row 10 ($r=9)  <- old sheet row 12 ($rowMap[9]=11)
#   sheet row 11 ($r=10) <- old sheet row 2  ($rowMap[10]=1)
#   sheet row 12 ($r=11) <- old sheet row 3  ($rowMap[11]=2)
$rowMap = @{
    1  = 3
    2  = 4
    3  = 5
    4  = 6
    5  = 7
    6  = 8
    7  = 9
    8  = 10
    9  = 11
    10 = 1
    11 = 2
}

$newArr = New-Object 'object[,]' $rowCount, $colCount

for ($r = 1; $r -le $rowCount; $r++) {
    $srcR = $rowMap[$r]
    for ($c = 1; $c -le $colCount; $c++) {
        $newArr[$r - 1, $c - 1] = $arr[$srcR, $c]
    }
}

# The Startdatum/Slutdatum columns (Y, AA) hold plain text dates such as
# "2023-06-13"; without forcing a text format first, Excel would
# reinterpret them as date serial numbers when the array is written
# back. Force text on just those two columns, write the values, then
# clear the temporary formatting so no stray number format/style is
# left behind (matching the original, unstyled inline-string cells).
$ws.Range("Y2:Y12").NumberFormat = "@"
$ws.Range("AA2:AA12").NumberFormat = "@"

$dstRange.Value = $newArr

$ws.Range("Y2:Y12").ClearFormats()
$ws.Range("AA2:AA12").ClearFormats()
